$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.340.21"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.604.00"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'592.47"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").Value = "'150.34"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").Value = "2.602.64"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").Value = "'27.27"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "3.077.13"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "67.166.21"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "2.605.75"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "'369.33"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "'11.03"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("E22").Value = "  -2.72%  "
$ws.Range("D23").Value = "'4.77"
$ws.Range("E23").Value = "  -3.67%  "
$ws.Range("E24").Value = "  -2.88%  "
$ws.Range("D25").Value = "'73.05"
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.730.55"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "'577.32"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").Value = "0.0₃0986"
$ws.Range("E31").Value = "  -5.57%  "
$ws.Range("E32").Value = "  -4.69%  "
$ws.Range("D33").Value = "'7.66"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "'157.52"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").Value = "'19.03"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E40").Value = "  -1.60%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("E43").Value = "  -3.61%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "'152.98"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").Value = "'21.31"
$ws.Range("E51").Value = "  +1.89%  "
